$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (row 4, 5, 6)
$ws.Range("A4").Value = 6664
$ws.Range("B4").Value = 4446

$ws.Range("A5").Value = 8000
$ws.Range("B5").Value = 6000

$ws.Range("A6").Value = 5800.9
$ws.Range("B6").Value = 6890.71

# Add new row 36 as text values "9000.0" (not numbers)
$ws.Range("A36").Value = "'9000.0"
$ws.Range("B36").Value = "'9000.0"
$ws.Range("A36:B36").ClearFormats()
